# Weekly update: insert a new price observation as row 62 (Rabanito,
# Vega Modelo de Temuco), pushing the existing rows 62-110 down to 63-111.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 62..110 down by one (Excel-style row insert).
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new observation.
$ws.Cells.Item(62, 1).Value2  = 10
$ws.Cells.Item(62, 2).Value   = "Vega Modelo de Temuco"
$ws.Cells.Item(62, 3).Value   = "La Araucanía"
$ws.Cells.Item(62, 4).Value2  = 45086
$ws.Cells.Item(62, 5).Value2  = 9
$ws.Cells.Item(62, 6).Value2  = 300000001
$ws.Cells.Item(62, 7).Value   = "Rabanito"
$ws.Cells.Item(62, 8).Value   = "Sin especificar"
$ws.Cells.Item(62, 9).Value   = "Primera"
$ws.Cells.Item(62, 10).Value2 = 40
$ws.Cells.Item(62, 11).Value2 = 7000
$ws.Cells.Item(62, 12).Value2 = 7000
$ws.Cells.Item(62, 13).Value2 = 7000
$ws.Cells.Item(62, 14).Value  = '$/docena de paquetes'
$ws.Cells.Item(62, 15).Value  = "Provincia de Cautín"
$ws.Cells.Item(62, 16).Value2 = 583
$ws.Cells.Item(62, 17).Value2 = 12
$ws.Cells.Item(62, 18).Value  = "Hortaliza"
